$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 841.48
$ws.Range("J129").Value = 931.0941
$ws.Range("L129").Value = 2793.2823
$ws.Range("N129").Value = -12793.2823

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4187.5835
$ws.Range("I61").Value = 4625.375
$ws.Range("K61").Value = 4625.375
$ws.Range("M61").Value = -4413.375

$ws.Range("H74").Value = 3369
$ws.Range("I74").Value = 3389.6667
$ws.Range("K74").Value = 3389.6667
$ws.Range("M74").Value = -2515.6667

$ws.Range("H77").Value = 3369
$ws.Range("I77").Value = 3389.6667
$ws.Range("K77").Value = 16948.3335
$ws.Range("M77").Value = -12580.3335

$ws.Range("H124").Value = 50400
$ws.Range("J124").Value = 50400
$ws.Range("L124").Value = 50400
$ws.Range("N124").Value = -60220

$ws.Range("H131").Value = 52197.3
$ws.Range("J131").Value = 52197.3
$ws.Range("L131").Value = 52197.3
$ws.Range("N131").Value = -62277.3

$ws.Range("H136").Value = 4187.5835
$ws.Range("I136").Value = 4625.375
$ws.Range("K136").Value = 13876.125
$ws.Range("M136").Value = -11326.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4134.641
$ws.Range("I31").Value = 1877.7059
$ws.Range("J31").Value = 5878.636
$ws.Range("K31").Value = 1877.7059
$ws.Range("L31").Value = 5878.636
$ws.Range("M31").Value = -1582.7059
$ws.Range("N31").Value = -6468.636

$ws.Range("H34").Value = 4134.641
$ws.Range("I34").Value = 1877.7059
$ws.Range("J34").Value = 5878.636
$ws.Range("K34").Value = 1877.7059
$ws.Range("L34").Value = 5878.636
$ws.Range("M34").Value = -1675.7059
$ws.Range("N34").Value = -6282.636

$ws.Range("H58").Value = 1836.9524
$ws.Range("I58").Value = 1573.8
$ws.Range("J58").Value = 2494.8333
$ws.Range("K58").Value = 1573.8
$ws.Range("L58").Value = 2494.8333
$ws.Range("M58").Value = -1370.8
$ws.Range("N58").Value = -2900.8333

$ws.Range("H99").Value = 2038.3125
$ws.Range("I99").Value = 1984.5385
$ws.Range("J99").Value = 2271.3333
$ws.Range("K99").Value = 1984.5385
$ws.Range("L99").Value = 2271.3333
$ws.Range("M99").Value = -486.5385000000001
$ws.Range("N99").Value = -5267.3333

$ws.Range("H126").Value = 2038.3125
$ws.Range("I126").Value = 1984.5385
$ws.Range("J126").Value = 2271.3333
$ws.Range("K126").Value = 5953.6155
$ws.Range("L126").Value = 6813.999899999999
$ws.Range("M126").Value = -3483.6155
$ws.Range("N126").Value = -11753.9999

$ws.Range("H136").Value = 1836.9524
$ws.Range("I136").Value = 1573.8
$ws.Range("J136").Value = 2494.8333
$ws.Range("K136").Value = 4721.4
$ws.Range("L136").Value = 7484.499899999999
$ws.Range("M136").Value = -2171.4
$ws.Range("N136").Value = -12584.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 435.17856
$ws.Range("I2").Value = 640.4375
$ws.Range("J2").Value = 161.5
$ws.Range("K2").Value = 3842.625
$ws.Range("L2").Value = 969
$ws.Range("M2").Value = -3729.625
$ws.Range("N2").Value = -1195

$ws.Range("H3").Value = 6061.769
$ws.Range("I3").Value = 3845
$ws.Range("K3").Value = 11535
$ws.Range("M3").Value = -11423

$ws.Range("H38").Value = 97.55
$ws.Range("J38").Value = 128
$ws.Range("L38").Value = 384
$ws.Range("N38").Value = -1078

$ws.Range("H122").Value = 916.73914
$ws.Range("I122").Value = 765.0909
$ws.Range("J122").Value = 1055.75
$ws.Range("K122").Value = 6885.8181
$ws.Range("L122").Value = 9501.75
$ws.Range("M122").Value = -4435.8181
$ws.Range("N122").Value = -14401.75

$ws.Range("H132").Value = 1771.4286
$ws.Range("I132").Value = 1420
$ws.Range("J132").Value = 2035
$ws.Range("K132").Value = 12780
$ws.Range("L132").Value = 18315
$ws.Range("M132").Value = -10250
$ws.Range("N132").Value = -23375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 350313.78
$ws.Range("I70").Value = 505175.25
$ws.Range("J70").Value = 6177.222
$ws.Range("K70").Value = 505175.25
$ws.Range("L70").Value = 6177.222
$ws.Range("M70").Value = -504905.25
$ws.Range("N70").Value = -6717.222

$ws.Range("H73").Value = 350313.78
$ws.Range("I73").Value = 505175.25
$ws.Range("J73").Value = 6177.222
$ws.Range("K73").Value = 505175.25
$ws.Range("L73").Value = 6177.222
$ws.Range("M73").Value = -504239.25
$ws.Range("N73").Value = -8049.222

$ws.Range("H109").Value = 15966.728
$ws.Range("J109").Value = 15966.728
$ws.Range("L109").Value = 15966.728
$ws.Range("N109").Value = -18046.728

$ws.Range("H123").Value = 837016
$ws.Range("J123").Value = 837016
$ws.Range("L123").Value = 837016
$ws.Range("N123").Value = -841916

$ws.Range("H132").Value = 3763.578
$ws.Range("I132").Value = 3160.1875
$ws.Range("K132").Value = 9480.5625
$ws.Range("M132").Value = -6950.5625

$ws.Range("H135").Value = 39826
$ws.Range("J135").Value = 39826
$ws.Range("L135").Value = 39826
$ws.Range("N135").Value = -49966

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H124").Value = 30429
$ws.Range("J124").Value = 30429
$ws.Range("L124").Value = 30429
$ws.Range("N124").Value = -40249

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3914.5715
$ws.Range("I62").Value = 3480.4
$ws.Range("K62").Value = 3480.4
$ws.Range("M62").Value = -2856.4

$ws.Range("H65").Value = 3914.5715
$ws.Range("I65").Value = 3480.4
$ws.Range("K65").Value = 17402
$ws.Range("M65").Value = -14282

$ws.Range("H68").Value = 49990.332
$ws.Range("J68").Value = 49990.332
$ws.Range("L68").Value = 49990.332
$ws.Range("N68").Value = -51612.332

$ws.Range("H69").Value = 20000
$ws.Range("J69").Value = 20000
$ws.Range("L69").Value = 20000
$ws.Range("N69").Value = -21498

$ws.Range("H71").Value = 49990.332
$ws.Range("J71").Value = 49990.332
$ws.Range("L71").Value = 149970.996
$ws.Range("N71").Value = -158082.996

$ws.Range("H72").Value = 20000
$ws.Range("J72").Value = 20000
$ws.Range("L72").Value = 60000
$ws.Range("N72").Value = -67488

$ws.Range("H125").Value = 49932.668
$ws.Range("J125").Value = 49932.668
$ws.Range("L125").Value = 49932.668
$ws.Range("N125").Value = -59772.668

$ws.Range("H131").Value = 59957
$ws.Range("J131").Value = 59957
$ws.Range("L131").Value = 59957
$ws.Range("N131").Value = -70037

$ws.Range("H132").Value = 4774.381
$ws.Range("I132").Value = 4663.9165
$ws.Range("J132").Value = 4921.6665
$ws.Range("K132").Value = 13991.7495
$ws.Range("L132").Value = 14764.9995
$ws.Range("M132").Value = -11461.7495
$ws.Range("N132").Value = -19824.9995

$ws.Range("H136").Value = 1985.62
$ws.Range("I136").Value = 1627.4412
$ws.Range("J136").Value = 2746.75
$ws.Range("K136").Value = 4882.3236
$ws.Range("L136").Value = 8240.25
$ws.Range("M136").Value = -2332.3236
$ws.Range("N136").Value = -13340.25
